# api for credit card repayment
# Applies the OOXML diff: renames the "repayment" endpoints to "repayments",
# adds a new "user wants to view repayments" (GET) row, shifts the
# statements / instalment-contract rows down by one row, and simplifies
# the instalment-contract(s) scenario text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, shifting rows 14-41 down to 15-42.
# This turns the old 2-row "repayment" block (rows 12-13) + blank (row 14)
# into a 3-row "repayments" block (rows 12-14) + blank (row 15), and keeps
# the table below (old rows 15-28) lined up at new rows 16-29.
$ws.Rows.Item(14).Insert()

# --- Repayment block (rows 12-14) ---
# Row 12: POST /users/{uid}/cards/{cid}/repayments  -> user wants to repay to a card
$ws.Range("B12").Value = "/users/{uid}/cards/{cid}/repayments"

# Row 14: GET /users/{uid}/cards/{cid}/repayments/{rid} -> user wants to view a repay to a card
# (populated before row 13's new text so new shared strings are appended in
# the same order as the source workbook)
$ws.Range("A14").Value = "user:{user_id}"
$ws.Range("B14").Value = "/users/{uid}/cards/{cid}/repayments/{rid}"
$ws.Range("C14").Value = "GET"
$ws.Range("D14").Value = "user wants to view a repay to a card"
$ws.Range("E14").Value = 200404

# Row 13: GET /users/{uid}/cards/{cid}/repayments -> user wants to view repayments (new row)
$ws.Range("A13").Value = "user:{user_id}"
$ws.Range("B13").Value = "/users/{uid}/cards/{cid}/repayments"
$ws.Range("C13").Value = "GET"
$ws.Range("D13").Value = "user wants to view repayments"
$ws.Range("E13").Value = 200404

# --- Simplify scenario text for instalment-contract(s) rows (new rows 26-27) ---
$ws.Range("D26").Value = "view instalment contracts"
$ws.Range("D27").Value = "view instalment contract"

# --- Update the saved selection to match the edited workbook ---
$ws.Range("D29").Select()
